# update 2x2 scalar test
#
# The "TwoxTwoScalar" sheet's results table had its two blocks of rows in
# the wrong order: the normalized-price block (PX.L/PX.L .. RA.L/PX.L) was
# above the raw quantity block (LX.L .. DY.L). This swaps the blocks back
# to the intended order (quantities first, then normalized prices) and
# fills in the previously-missing benchmark (col B) and third scenario
# (col E) values for the quantity rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TwoxTwoScalar")

# Move the quantity block (rows 17:22 -> LX.L, LY.L, KX.L, KY.L, DX.L, DY.L)
# out of the way, then shift the price-ratio block (rows 11:16) down to
# 17:22, then drop the quantity block into its new home at 11:16. Using a
# scratch area (row 30) avoids the two blocks overlapping mid-move.
$ws.Range("A17:E22").Cut($ws.Range("A30"))
$ws.Range("A11:E16").Cut($ws.Range("A17"))
$ws.Range("A30:E35").Cut($ws.Range("A11"))

# Remove the now-empty scratch rows so the sheet dimension shrinks back.
$ws.Rows("23:35").Delete()

# Fill in the benchmark (B) and third-scenario (E) values for the
# quantity rows, which previously only had columns C and D populated.
$ws.Range("B11").Value = 50
$ws.Range("E11").Value = 52.440442408500388

$ws.Range("B12").Value = 20
$ws.Range("E12").Value = 21.177057058432069

$ws.Range("B13").Value = 50
$ws.Range("E13").Value = 47.673129462283626

$ws.Range("B14").Value = 30
$ws.Range("E14").Value = 28.877805079687114

$ws.Range("B15").Value = 100
$ws.Range("E15").Value = 100.31820580257069

$ws.Range("B16").Value = 50
$ws.Range("E16").Value = 49.683306602973595

# The active tab moves from "TwoxTwoAlg" back to "TwoxTwoScalar".
$ws.Activate()
